$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Tue Jul 11 12:55:13 EDT 2023"
$ws.Range("B3").Value = "Tue Jul 11 12:55:23 EDT 2023"
$ws.Range("B4").Value = "Tue Jul 11 12:55:33 EDT 2023"
